$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# The "Periodo Mora" column (E16:E23) is re-populated with the same set of
# periods but in descending order (newest period first), and the "Valor
# Mora" figure that used to be tied to period 2103 (row 23) now travels
# with whichever row holds that period (row 16) - i.e. F16/F23 swap.
$ws.Range("E16").Value = "2103"
$ws.Range("E17").Value = "2102"
$ws.Range("E18").Value = "2101"
$ws.Range("E19").Value = "2012"
$ws.Range("E20").Value = "2011"
$ws.Range("E21").Value = "2010"
$ws.Range("E22").Value = "2009"
$ws.Range("E23").Value = "2008"

$ws.Range("F16").Value = 26919
$ws.Range("F23").Value = 35112
